$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-parsed by Excel as numbers (single decimal point numerics),
# so they remain plain text as in the source data.
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "26.092.00"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.648.00"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "216.92"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "0.0641"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").Value = "19.66"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.875.48"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "1.679.13"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.30"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "0.545"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "26.218.66"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "195.90"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "9.95"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "6.24"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "0.131"
$ws.Range("E25").Value = "  +4.14%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").Value = "144.25"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "6.92"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "15.59"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "1.25"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").Value = "3.29"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "2.46"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "0.908"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.136.32"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.543"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.46"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0157"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.51"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "99.53"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.799"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.782.73"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0117"
$ws.Range("E45").Value = "  +4.74%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "56.77"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0527"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.46"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.76"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.417"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0960"
$ws.Range("E51").Value = "  -0.08%  "
